{"js": "// Adds a new \"ex\" (example) hyperlink run after the existing \"vid\"/\"vid1\"/\"vid2\"\n// video-hyperlink run(s) in several CSS3 topic bullets (FILTER, TRANSFORM,\n// TRANSITION, ANIMATION, GALERI FOTO + LIGHTBOX), mirroring the \"vid ... ex\"\n// pattern already used by every other bullet in the document.\n//\n// For each target paragraph we:\n//   1. locate it by its distinctive leading text,\n//   2. insert \"  \" + the new link's display text at the end of the paragraph,\n//   3. search for that freshly-inserted text inside the paragraph, and\n//   4. set the `.hyperlink` property on the found range, which turns the run\n//      into a proper `w:hyperlink` (styled with the built-in Hyperlink\n//      character style), exactly like the pre-existing \"vid\" links.\n\n// Plan: paragraph (matched by the text it starts with) -> list of links to\n// append, in left-to-right order, after its current content.\nconst plan = [\n  { startsWith: \"FILTER\", links: [\n    { text: \"ex\", url: \"ex/filter.htm\" },\n  ] },\n  { startsWith: \"TRANSFORM\", links: [\n    { text: \"ex\", url: \"ex/transform.htm\" },\n  ] },\n  { startsWith: \"TRANSITION\", links: [\n    { text: \"ex1\", url: \"ex/transition1.htm\" },\n    { text: \"ex2\", url: \"ex/transition2.htm\" },\n  ] },\n  { startsWith: \"ANIMATION\", links: [\n    { text: \"ex\", url: \"ex/animation.htm\" },\n  ] },\n  { startsWith: \"GALERI\", links: [\n    { text: \"ex1\", url: \"ex/lightbox1.htm\" },\n    { text: \"ex2\", url: \"ex/lightbox2.htm\" },\n  ] },\n];\n\n// Resolve each target paragraph's current index by scanning the body once.\nasync function findParagraphIndex(startsWith) {\n  const paras = context.document.body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  for (const p of paras.items) {\n    p.load(\"text\");\n  }\n  await context.sync();\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text.trim().indexOf(startsWith) === 0) {\n      return i;\n    }\n  }\n  throw new Error(\"Paragraph starting with '\" + startsWith + \"' not found\");\n}\n\nfor (const entry of plan) {\n  const paraIndex = await findParagraphIndex(entry.startsWith);\n\n  for (const link of entry.links) {\n    // 1) Append \"  \" + display text as plain text at the end of the paragraph.\n    const parasBefore = context.document.body.paragraphs;\n    parasBefore.load(\"items\");\n    await context.sync();\n    const para = parasBefore.items[paraIndex];\n    const endRange = para.getRange(\"End\");\n    endRange.insertText(\"  \" + link.text, \"End\");\n    await context.sync();\n\n    // 2) Re-fetch the paragraph and search for the text we just inserted so we\n    // get a live Range over exactly that run (insertText's returned range\n    // isn't reliably reusable across this host's sync boundary).\n    const parasAfter = context.document.body.paragraphs;\n    parasAfter.load(\"items\");\n    await context.sync();\n    const paraAfter = parasAfter.items[paraIndex];\n    const matches = paraAfter.search(link.text, { matchCase: true, matchWholeWord: false });\n    matches.load(\"items\");\n    await context.sync();\n\n    // The text we just appended is the last occurrence in the paragraph.\n    const target = matches.items[matches.items.length - 1];\n\n    // 3) Turning the range into a hyperlink produces a `w:hyperlink` wrapping\n    // a run styled with the built-in \"Hyperlink\" character style - matching\n    // the existing \"vid\" links exactly.\n    target.hyperlink = link.url;\n    await context.sync();\n  }\n}\n", "ps1": "# Adds a new \"ex\" (example) hyperlink run after the existing \"vid\"/\"vid1\"/\"vid2\"\n# video-hyperlink run(s) in several CSS3 topic bullets (FILTER, TRANSFORM,\n# TRANSITION, ANIMATION, GALERI FOTO + LIGHTBOX), mirroring the \"vid ... ex\"\n# pattern already used by every other bullet in the document.\n\n$d = $word.ActiveDocument\n\nfunction Add-ExLink {\n    param(\n        [string]$SearchText,   # distinctive text used to locate the target paragraph\n        [string]$LinkText,     # display text of the new hyperlink (\"ex\", \"ex1\", \"ex2\", ...)\n        [string]$Url           # hyperlink target\n    )\n\n    # Locate the paragraph that contains $SearchText.\n    $rng = $d.Content\n    $found = $rng.Find.Execute($SearchText)\n    if (-not $found) {\n        throw \"Paragraph containing '$SearchText' not found\"\n    }\n    $para = $rng.Paragraphs(1)\n    $pRange = $para.Range\n\n    # Insert \"  \" + a one-character placeholder right before the paragraph\n    # mark. Word COM's Hyperlinks.Add ignores the position of a *collapsed*\n    # (zero-length) range and always inserts at the start of the paragraph,\n    # so we first materialize a one-character range at the correct spot and\n    # then hand that non-empty range to Hyperlinks.Add, which replaces it\n    # in place with the hyperlink's display text.\n    $insertPos = $pRange.End - 1\n    $ins = $d.Range($insertPos, $insertPos)\n    $placeholder = \"  \" + [char]1\n    $ins.InsertAfter($placeholder)\n\n    $xPos = $insertPos + 2\n    $xRange = $d.Range($xPos, $xPos + 1)\n    $d.Hyperlinks.Add($xRange, $Url, \"\", \"\", $LinkText) | Out-Null\n}\n\nAdd-ExLink \"FILTER\" \"ex\" \"ex/filter.htm\"\nAdd-ExLink \"TRANSFORM\" \"ex\" \"ex/transform.htm\"\nAdd-ExLink \"TRANSITION\" \"ex1\" \"ex/transition1.htm\"\nAdd-ExLink \"TRANSITION\" \"ex2\" \"ex/transition2.htm\"\nAdd-ExLink \"ANIMATION\" \"ex\" \"ex/animation.htm\"\nAdd-ExLink \"GALERI\" \"ex1\" \"ex/lightbox1.htm\"\nAdd-ExLink \"GALERI\" \"ex2\" \"ex/lightbox2.htm\"\n"}
